$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value in D1 (plain number, not a shared string)
$ws.Range("D1").Value = 500081

# Update the selection to match the authored state (C2 selected)
$ws.Range("C2").Select()
